$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("A2").Style

$ws.Range("D2").Value = "65.998.57"
$ws.Range("E2").Value = "  +6.69%  "

$ws.Range("D3").Value = "3.017.46"
$ws.Range("E3").Value = "  +3.85%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'585.06"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +2.72%  "

$ws.Range("D6").Value = "'161.65"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +12.52%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.014.39"
$ws.Range("E8").Value = "  +3.86%  "

$ws.Range("D9").Value = "'0.517"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  +3.51%  "

$ws.Range("D10").Value = "'6.77"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("E11").Value = "  +6.54%  "

$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +7.41%  "

$ws.Range("D13").Value = "'0.0000254"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +9.51%  "

$ws.Range("D14").Value = "'34.68"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +7.46%  "

$ws.Range("E15").Value = "  -0.58%  "

$ws.Range("D16").Value = "65.952.85"
$ws.Range("E16").Value = "  +6.67%  "

$ws.Range("D17").Value = "3.517.22"
$ws.Range("E17").Value = "  +3.84%  "

$ws.Range("D18").Value = "'6.96"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  +7.56%  "

$ws.Range("D19").Value = "3.013.66"
$ws.Range("E19").Value = "  +3.61%  "

$ws.Range("D20").Value = "'457.64"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +6.73%  "

$ws.Range("D21").Value = "'13.99"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +7.91%  "

$ws.Range("E22").Value = "  +6.17%  "

$ws.Range("D23").Value = "'7.39"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +7.82%  "

$ws.Range("D24").Value = "'82.42"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +4.78%  "

$ws.Range("D25").Value = "'2.28"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  +12.89%  "

$ws.Range("D26").Value = "'12.40"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +3.51%  "

$ws.Range("D27").Value = "'10.65"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +4.75%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'8.05"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  +15.70%  "

$ws.Range("E30").Value = "  +16.17%  "

$ws.Range("E31").Value = "  -6.04%  "

$ws.Range("E32").Value = "  +4.08%  "

$ws.Range("D33").Value = "'27.14"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +6.35%  "

$ws.Range("E34").Value = "  +4.33%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").Value = "'0.997"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +4.29%  "

$ws.Range("D37").Value = "'5.84"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  +8.79%  "

$ws.Range("E38").Value = "  +14.72%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.00"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +3.71%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'49.88"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  +2.24%  "

$ws.Range("D41").Value = "'0.310"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  +16.76%  "

$ws.Range("E42").Value = "  +7.11%  "

$ws.Range("D43").Value = "'43.50"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +6.15%  "

$ws.Range("D44").Value = "'8.47"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +4.03%  "

$ws.Range("D45").Value = "'390.22"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +12.30%  "

$ws.Range("D46").Value = "2.801.48"
$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("D47").Value = "'0.0356"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +5.58%  "

$ws.Range("D48").Value = "'134.73"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +1.63%  "

$ws.Range("D50").Value = "'23.67"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  +9.79%  "

$ws.Range("E51").Value = "  +4.51%  "
